$d = $word.ActiveDocument

# Update the date heading (first paragraph).
$d.Content.Find.Execute("2024-01-20 Saturday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-01-21 Sunday", 2)

# Update the division problems in the table. Addressed by (row, column) so
# that values being reused elsewhere in the grid (e.g. "75÷8=") can't be
# double-matched the way a sequential global Find/Replace could.
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "97÷9=" },
    @{ Row = 1;  Col = 2; Text = "19÷6=" },
    @{ Row = 1;  Col = 3; Text = "46÷9=" },
    @{ Row = 1;  Col = 4; Text = "35÷5=" },
    @{ Row = 1;  Col = 5; Text = "75÷9=" },

    @{ Row = 5;  Col = 1; Text = "61÷5=" },
    @{ Row = 5;  Col = 2; Text = "25÷3=" },
    @{ Row = 5;  Col = 3; Text = "75÷8=" },
    @{ Row = 5;  Col = 4; Text = "91÷6=" },
    @{ Row = 5;  Col = 5; Text = "23÷2=" },

    @{ Row = 9;  Col = 1; Text = "82÷3=" },
    @{ Row = 9;  Col = 2; Text = "68÷2=" },
    @{ Row = 9;  Col = 3; Text = "96÷5=" },
    @{ Row = 9;  Col = 4; Text = "36÷8=" },
    @{ Row = 9;  Col = 5; Text = "45÷8=" },

    @{ Row = 13; Col = 1; Text = "64÷8=" },
    @{ Row = 13; Col = 2; Text = "49÷7=" },
    @{ Row = 13; Col = 3; Text = "21÷3=" },
    @{ Row = 13; Col = 4; Text = "64÷2=" },
    @{ Row = 13; Col = 5; Text = "13÷8=" },

    @{ Row = 17; Col = 1; Text = "43÷4=" },
    @{ Row = 17; Col = 2; Text = "50÷2=" },
    @{ Row = 17; Col = 3; Text = "87÷5=" },
    @{ Row = 17; Col = 4; Text = "94÷9=" },
    @{ Row = 17; Col = 5; Text = "77÷6=" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $r = $cell.Range
    $r.End = $r.End - 1
    $r.Text = $u.Text
}
